# Atualizei dados da bibi
# Insert a new daily sales record at the top of the "06/2025" block (row 10),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 10 (shifts rows 10:101 down to 11:102)
$ws.Rows.Item(10).Insert()

# Populate the new row with the inserted daily record
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 23856.87
$ws.Cells.Item(10, 3).Value = 6
$ws.Cells.Item(10, 4).Value = 2025
$ws.Cells.Item(10, 5).Value = "06/2025"
